# Apply changes to the "structure" sheet:
# - personId row: size 6 -> 3, format normal -> long, requiredforupdate yes -> no
# - firstNameActual row: format uppercase -> lowercasealphanumeric
# - remove the data validation list on D2:D5 (format column)
# - remove right border on D2:D5
# - move selection to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("structure")

# --- Content changes ---
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "long"
$ws.Range("G2").Value = "no"
$ws.Range("D3").Value = "lowercasealphanumeric"

# --- Remove data validation on D2:D5 ---
$ws.Range("D2:D5").Validation.Delete()

# --- Remove right border from D2:D5 (keep left/top/bottom) ---
$ws.Range("D2:D5").Borders(10).LineStyle = -4142

# --- Update selection ---
$ws.Range("D11").Select()
